$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$shape.TextFrame.TextRange.Text = "Let's exercise more!"
